# Updated cryptos list on Fri Jun 28 19:22:25 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the cell to hold an exact text value (avoid Excel's automatic
    # number/date coercion for numeric-looking strings like "1.00" or
    # "0.0000170"), then restore the default "Normal" style so no stray
    # number-format style is left attached to the cell.
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "60.422.69"
$ws.Range("E2").Value = "  -1.86%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.371.48"
$ws.Range("E3").Value = "  -2.28%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.10%  "

# Row 5 - BNB
Set-TextValue "D5" "570.42"
$ws.Range("E5").Value = "  -1.55%  "

# Row 6 - Solana
Set-TextValue "D6" "141.34"

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.06%  "

# Row 8 - LidoStakedEther
Set-TextValue "D8" "3.373.09"
$ws.Range("E8").Value = "  -2.29%  "

# Row 9 - XRP
Set-TextValue "D9" "0.473"
$ws.Range("E9").Value = "  +0.00%  "

# Row 10 - Toncoin
Set-TextValue "D10" "7.47"
$ws.Range("E10").Value = "  -3.87%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -0.65%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  +0.14%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue "D13" "3.958.68"
$ws.Range("E13").Value = "  -2.06%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  +1.13%  "

# Row 15 - Avalanche
Set-TextValue "D15" "27.96"
$ws.Range("E15").Value = "  -0.11%  "

# Row 16 - swapped to ShibaInu
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D16" "0.0000170"
$ws.Range("E16").Value = "  -2.58%  "

# Row 17 - swapped to WrappedEther
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D17" "3.380.93"
$ws.Range("E17").Value = "  -1.68%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "60.707.52"
$ws.Range("E18").Value = "  -1.62%  "

# Row 19 - Polkadot
Set-TextValue "D19" "6.25"
$ws.Range("E19").Value = "  -1.16%  "

# Row 20 - Chainlink
Set-TextValue "D20" "14.01"
$ws.Range("E20").Value = "  -2.14%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  -4.14%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "388.20"
$ws.Range("E22").Value = "  +0.58%  "

# Row 23 - Polygon
Set-TextValue "D23" "0.558"
$ws.Range("E23").Value = "  -2.06%  "

# Row 24 - Litecoin
Set-TextValue "D24" "73.18"
$ws.Range("E24").Value = "  +0.71%  "

# Row 25 - Dai
Set-TextValue "D25" "0.998"
$ws.Range("E25").Value = "  -0.25%  "

# Row 26 - PEPE
$ws.Range("E26").Value = "  -4.49%  "

# Row 27 - WrappedeETH
Set-TextValue "D27" "3.521.48"
$ws.Range("E27").Value = "  -2.12%  "

# Row 28 - Kaspa
$ws.Range("E28").Value = "  -0.65%  "

# Row 29 - Binance-PegBSC-USD
Set-TextValue "D29" "1.00"
$ws.Range("E29").Value = "  +0.03%  "

# Row 30 - RenderToken
Set-TextValue "D30" "7.39"
$ws.Range("E30").Value = "  -5.50%  "

# Row 31 - InternetComputer(DFINITY)
Set-TextValue "D31" "8.05"
$ws.Range("E31").Value = "  -2.34%  "

# Row 32 - PancakeSwap
$ws.Range("E32").Value = "  -0.68%  "

# Row 33 - Fetch.AI
$ws.Range("E33").Value = "  -4.09%  "

# Row 34 - USDe
$ws.Range("E34").Value = "  -0.03%  "

# Row 35 - EthereumClassic
Set-TextValue "D35" "23.67"
$ws.Range("E35").Value = "  -1.31%  "

# Row 37 - RenzoRestakedETH
Set-TextValue "D37" "3.410.09"
$ws.Range("E37").Value = "  -1.84%  "

# Row 38 - Monero
Set-TextValue "D38" "166.93"
$ws.Range("E38").Value = "  +0.42%  "

# Row 39 - NEARProtocol
Set-TextValue "D39" "4.98"
$ws.Range("E39").Value = "  -4.44%  "

# Row 40 - ImmutableX
Set-TextValue "D40" "1.50"
$ws.Range("E40").Value = "  -3.88%  "

# Row 41 - Hedera
$ws.Range("E41").Value = "  -1.60%  "

# Row 42 - EnergySwap
Set-TextValue "D42" "26.70"
$ws.Range("E42").Value = "  +2.32%  "

# Row 43 - swapped to FirstDigitalUSD
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D43" "1.00"
$ws.Range("E43").Value = "  +0.21%  "

# Row 44 - swapped to Mantle
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D44" "0.781"
$ws.Range("E44").Value = "  -1.95%  "

# Row 45 - Filecoin
Set-TextValue "D45" "4.44"
$ws.Range("E45").Value = "  -1.24%  "

# Row 46 - OKB
Set-TextValue "D46" "41.71"
$ws.Range("E46").Value = "  -1.49%  "

# Row 47 - Stacks
$ws.Range("E47").Value = "  -2.13%  "

# Row 48 - Maker
Set-TextValue "D48" "2.532.55"
$ws.Range("E48").Value = "  -2.38%  "

# Row 49 - ONDO
$ws.Range("E49").Value = "  -3.69%  "

# Row 50 - Cosmos
Set-TextValue "D50" "6.80"
$ws.Range("E50").Value = "  -2.44%  "

# Row 51 - InjectiveProtocol
Set-TextValue "D51" "22.84"
$ws.Range("E51").Value = "  -1.96%  "
